$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell explicitly as text, preserving original formatting/style
# (ClearFormats after a forced "@" text NumberFormat keeps the stored value as a
# string while restoring the default/unstyled cell format, matching the source file
# where these cells carry no explicit style.)
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("B37") "ImmutableX"
Set-TextValue $ws.Range("B38") "Maker"
Set-TextValue $ws.Range("B43") "BabyDogeCoin"
Set-TextValue $ws.Range("B44") "Quant"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D2") "25.945.44"
Set-TextValue $ws.Range("D3") "1.646.84"
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("D5") "216.02"
Set-TextValue $ws.Range("D6") "0.5106"
Set-TextValue $ws.Range("D7") "1.006"
Set-TextValue $ws.Range("D8") "0.2585"
Set-TextValue $ws.Range("D9") "0.06431"
Set-TextValue $ws.Range("D10") "19.70"
Set-TextValue $ws.Range("D11") "0.07775"
Set-TextValue $ws.Range("D13") "1.652.76"
Set-TextValue $ws.Range("D14") "0.5478"
Set-TextValue $ws.Range("D15") "0.0₅7904"
Set-TextValue $ws.Range("D16") "64.92"
Set-TextValue $ws.Range("D17") "26.017.06"
Set-TextValue $ws.Range("D18") "1.006"
Set-TextValue $ws.Range("D19") "198.78"
Set-TextValue $ws.Range("D20") "4.468"
Set-TextValue $ws.Range("D22") "6.081"
Set-TextValue $ws.Range("D24") "1.864"
Set-TextValue $ws.Range("D25") "140.38"
Set-TextValue $ws.Range("D27") "6.915"
Set-TextValue $ws.Range("D29") "1.242"
Set-TextValue $ws.Range("D30") "0.05038"
Set-TextValue $ws.Range("D31") "3.293"
Set-TextValue $ws.Range("D32") "3.207"
Set-TextValue $ws.Range("D33") "1.546"
Set-TextValue $ws.Range("D34") "2.363"
Set-TextValue $ws.Range("D35") "0.8952"
Set-TextValue $ws.Range("D36") "2.592"
Set-TextValue $ws.Range("D37") "0.5553"
Set-TextValue $ws.Range("D38") "1.135.36"
Set-TextValue $ws.Range("D39") "0.01577"
Set-TextValue $ws.Range("D40") "1.006"
Set-TextValue $ws.Range("D41") "5.677"
Set-TextValue $ws.Range("D42") "0.8184"
Set-TextValue $ws.Range("D43") "0.0₈126"
Set-TextValue $ws.Range("D44") "100.02"
Set-TextValue $ws.Range("D45") "1.787.04"
Set-TextValue $ws.Range("D46") "0.4536"
Set-TextValue $ws.Range("D47") "55.34"
Set-TextValue $ws.Range("D49") "0.05094"
Set-TextValue $ws.Range("D50") "0.09591"
Set-TextValue $ws.Range("D51") "1.006"
Set-TextValue $ws.Range("E2") "  +0.44%  "
Set-TextValue $ws.Range("E3") "  +0.81%  "
Set-TextValue $ws.Range("E4") "  +0.48%  "
Set-TextValue $ws.Range("E5") "  +0.62%  "
Set-TextValue $ws.Range("E6") "  +1.73%  "
Set-TextValue $ws.Range("E7") "  +0.39%  "
Set-TextValue $ws.Range("E8") "  +0.85%  "
Set-TextValue $ws.Range("E9") "  +0.68%  "
Set-TextValue $ws.Range("E10") "  +0.75%  "
Set-TextValue $ws.Range("E11") "  +1.18%  "
Set-TextValue $ws.Range("E12") "  +2.02%  "
Set-TextValue $ws.Range("E13") "  +1.21%  "
Set-TextValue $ws.Range("E14") "  +1.03%  "
Set-TextValue $ws.Range("E15") "  -0.09%  "
Set-TextValue $ws.Range("E16") "  +2.36%  "
Set-TextValue $ws.Range("E17") "  +0.67%  "
Set-TextValue $ws.Range("E18") "  +0.30%  "
Set-TextValue $ws.Range("E19") "  -1.44%  "
Set-TextValue $ws.Range("E20") "  +3.49%  "
Set-TextValue $ws.Range("E21") "  +1.35%  "
Set-TextValue $ws.Range("E22") "  +1.92%  "
Set-TextValue $ws.Range("E23") "  +0.42%  "
Set-TextValue $ws.Range("E24") "  -2.33%  "
Set-TextValue $ws.Range("E25") "  -0.60%  "
Set-TextValue $ws.Range("E26") "  +1.46%  "
Set-TextValue $ws.Range("E27") "  +3.35%  "
Set-TextValue $ws.Range("E28") "  +0.80%  "
Set-TextValue $ws.Range("E29") "  +0.33%  "
Set-TextValue $ws.Range("E30") "  +1.32%  "
Set-TextValue $ws.Range("E31") "  +1.21%  "
Set-TextValue $ws.Range("E32") "  +1.18%  "
Set-TextValue $ws.Range("E33") "  +0.81%  "
Set-TextValue $ws.Range("E34") "  -0.09%  "
Set-TextValue $ws.Range("E35") "  +0.52%  "
Set-TextValue $ws.Range("E36") "  -1.05%  "
Set-TextValue $ws.Range("E37") "  -0.65%  "
Set-TextValue $ws.Range("E38") "  -2.84%  "
Set-TextValue $ws.Range("E39") "  +1.23%  "
Set-TextValue $ws.Range("E40") "  +0.43%  "
Set-TextValue $ws.Range("E41") "  -0.02%  "
Set-TextValue $ws.Range("E42") "  +1.43%  "
Set-TextValue $ws.Range("E43") "  +10.04%  "
Set-TextValue $ws.Range("E44") "  +0.74%  "
Set-TextValue $ws.Range("E45") "  +0.95%  "
Set-TextValue $ws.Range("E46") "  +0.46%  "
Set-TextValue $ws.Range("E47") "  +1.39%  "
Set-TextValue $ws.Range("E48") "  +0.44%  "
Set-TextValue $ws.Range("E49") "  +0.32%  "
Set-TextValue $ws.Range("E51") "  +0.28%  "
